# Apply the cryptocurrency price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value, and whether Excel would
# otherwise misinterpret the text as a number (needs to be forced to text).
$updates = @(
    @{ Cell = "D2"; Value = "43.086.35"; ForceText = $false },
    @{ Cell = "E2"; Value = "  +1.44%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "2.375.31"; ForceText = $false },
    @{ Cell = "E3"; Value = "  +6.68%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  -0.40%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "323.30"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +9.70%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "103.59"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -7.28%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "0.641"; ForceText = $true },
    @{ Cell = "E7"; Value = "  +2.43%  "; ForceText = $false },
    @{ Cell = "E8"; Value = "  -0.04%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "0.656"; ForceText = $true },
    @{ Cell = "E9"; Value = "  +10.03%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "41.40"; ForceText = $true },
    @{ Cell = "E10"; Value = "  -4.42%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "0.0936"; ForceText = $true },
    @{ Cell = "E11"; Value = "  +1.85%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "8.54"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -1.63%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "1.02"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -2.20%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "16.86"; ForceText = $true },
    @{ Cell = "E14"; Value = "  +12.86%  "; ForceText = $false },
    @{ Cell = "E15"; Value = "  +2.12%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "2.736.83"; ForceText = $false },
    @{ Cell = "E16"; Value = "  +6.82%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "2.379.79"; ForceText = $false },
    @{ Cell = "E17"; Value = "  +8.03%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "43.108.23"; ForceText = $false },
    @{ Cell = "E18"; Value = "  +1.54%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "7.85"; ForceText = $true },
    @{ Cell = "E19"; Value = "  +9.35%  "; ForceText = $false },
    @{ Cell = "E20"; Value = "  +2.13%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "76.46"; ForceText = $true },
    @{ Cell = "E21"; Value = "  +4.10%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "276.00"; ForceText = $true },
    @{ Cell = "E22"; Value = "  +15.26%  "; ForceText = $false },
    @{ Cell = "E23"; Value = "  +0.24%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "2.39"; ForceText = $true },
    @{ Cell = "E24"; Value = "  +1.58%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "9.50"; ForceText = $true },
    @{ Cell = "E25"; Value = "  +7.34%  "; ForceText = $false },
    @{ Cell = "E26"; Value = "  +0.09%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "23.09"; ForceText = $true },
    @{ Cell = "E28"; Value = "  +6.94%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "37.91"; ForceText = $true },
    @{ Cell = "E29"; Value = "  +2.42%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "175.56"; ForceText = $true },
    @{ Cell = "E30"; Value = "  +0.23%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "2.16"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -1.35%  "; ForceText = $false },
    @{ Cell = "E32"; Value = "  +2.11%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "0.0917"; ForceText = $true },
    @{ Cell = "E33"; Value = "  +4.84%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "5.83"; ForceText = $true },
    @{ Cell = "E34"; Value = "  +2.84%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "0.132"; ForceText = $true },
    @{ Cell = "E35"; Value = "  +4.88%  "; ForceText = $false },
    @{ Cell = "D36"; Value = "4.85"; ForceText = $true },
    @{ Cell = "E36"; Value = "  -1.37%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "4.18"; ForceText = $true },
    @{ Cell = "E37"; Value = "  +0.43%  "; ForceText = $false },
    @{ Cell = "E38"; Value = "  -2.12%  "; ForceText = $false },
    @{ Cell = "E39"; Value = "  +1.45%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "2.80"; ForceText = $true },
    @{ Cell = "E40"; Value = "  +17.58%  "; ForceText = $false },
    @{ Cell = "E41"; Value = "  +21.05%  "; ForceText = $false },
    @{ Cell = "B42"; Value = "Aave"; ForceText = $false },
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; ForceText = $false },
    @{ Cell = "D42"; Value = "123.60"; ForceText = $true },
    @{ Cell = "E42"; Value = "  +21.89%  "; ForceText = $false },
    @{ Cell = "B43"; Value = "Algorand"; ForceText = $false },
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; ForceText = $false },
    @{ Cell = "D43"; Value = "0.229"; ForceText = $true },
    @{ Cell = "E43"; Value = "  +1.00%  "; ForceText = $false },
    @{ Cell = "B44"; Value = "MultiversX"; ForceText = $false },
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"; ForceText = $false },
    @{ Cell = "D44"; Value = "69.15"; ForceText = $true },
    @{ Cell = "E44"; Value = "  -2.78%  "; ForceText = $false },
    @{ Cell = "B45"; Value = "FirstDigitalUSD"; ForceText = $false },
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; ForceText = $false },
    @{ Cell = "D45"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E45"; Value = "  +0.16%  "; ForceText = $false },
    @{ Cell = "B46"; Value = "BitcoinSV"; ForceText = $false },
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"; ForceText = $false },
    @{ Cell = "D46"; Value = "92.91"; ForceText = $true },
    @{ Cell = "E46"; Value = "  +58.21%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "12.35"; ForceText = $true },
    @{ Cell = "E47"; Value = "  +0.99%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "9.48"; ForceText = $true },
    @{ Cell = "E48"; Value = "  +11.86%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "5.57"; ForceText = $true },
    @{ Cell = "E49"; Value = "  +2.91%  "; ForceText = $false },
    @{ Cell = "E50"; Value = "  +0.97%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "1.597.45"; ForceText = $false },
    @{ Cell = "E51"; Value = "  +11.75%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Temporarily format as text so strings like "0.641" or "323.30"
        # are kept verbatim instead of being coerced into numbers, then
        # clear the (now-unneeded) explicit formatting again.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
